$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new record as row 16, pushing the existing rows 16..79 down to 17..80.
# (Excel copies the row-above's formatting for the newly inserted row, matching
# the style of the surrounding date column.)
$ws.Rows("16:16").Insert()

# Populate the newly inserted row 16 with the new "Membrillo" market record.
$ws.Cells.Item(16, 1).Value = 6
$ws.Cells.Item(16, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(16, 3).Value = "Metropolitana"
$ws.Cells.Item(16, 4).Value = 44659
$ws.Cells.Item(16, 5).Value = 13
$ws.Cells.Item(16, 6).Value = "Fruta"
$ws.Cells.Item(16, 7).Value = 100104
$ws.Cells.Item(16, 8).Value = "Frutos de pepita"
$ws.Cells.Item(16, 9).Value = 100104003
$ws.Cells.Item(16, 10).Value = "Membrillo"
$ws.Cells.Item(16, 11).Value = "Champion"
$ws.Cells.Item(16, 12).Value = "Primera"
$ws.Cells.Item(16, 13).Value = 24
$ws.Cells.Item(16, 14).Value = 240000
$ws.Cells.Item(16, 15).Value = 250000
$ws.Cells.Item(16, 16).Value = 247083
$ws.Cells.Item(16, 17).Value = "$/bins (450 kilos)"
$ws.Cells.Item(16, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(16, 19).Value = 549
$ws.Cells.Item(16, 20).Value = 450
